$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the geocoded address strings for the new orders
$ws.Range("G1").Value = "Petisah, Kota Medan, Sumatera Utara, Sumatera, 20256, Indonesia"
$ws.Range("G3").Value = "Sukabumi II, Kota Medan, Sumatera Utara, Sumatera, 20222, Indonesia"

# Update the active selection to match the saved workbook state
$ws.Range("G4").Select()
